# Update "paises.xlsx" (sheet "Pais") with the latest COVID-19 country
# snapshot and refresh the "Provincias España" (sub-country breakdown)
# figures, per commit "Update countries & provincias Spain".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# ---------------------------------------------------------------------
# 1) Timestamp banner in row 1
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 25 de Septiembre de 2020 a las 23:10"

# ---------------------------------------------------------------------
# 2) Refreshed per-country counters (Casos totales, Nuevos casos,
#    Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
#    Columns: A Pais | B Casos totales | C Nuevos casos | D Casos activos
#             E Recuperados | F Casos criticos | G Muertes hoy | H Muertes
# ---------------------------------------------------------------------

function Set-Row($r, $b, $c, $d, $e, $f, $g, $h) {
    if ($b -ne $null) { $ws.Cells.Item($r, 2).Value = $b }
    if ($c -ne $null) { $ws.Cells.Item($r, 3).Value = $c }
    if ($d -ne $null) { $ws.Cells.Item($r, 4).Value = $d }
    if ($e -ne $null) { $ws.Cells.Item($r, 5).Value = $e }
    if ($f -ne $null) { $ws.Cells.Item($r, 6).Value = $f }
    if ($g -ne $null) { $ws.Cells.Item($r, 7).Value = $g }
    if ($h -ne $null) { $ws.Cells.Item($r, 8).Value = $h }
}

# Row 4 - Estados Unidos
Set-Row 4 7222931 37460 4464431 2550383 $null 597 208117
# Row 5 - India
Set-Row 5 $null $null $null 961998 $null 1093 93410
# Row 6 - Brasil
Set-Row 6 4689613 29704 $null 525287 $null 654 140537
# Row 25 - Alemania
Set-Row 25 283706 2361 $null 24676 $null $null $null
# Row 29 - Canada
Set-Row 29 150194 1100 129318 11621 $null 6 9255
# Row 31 - Ecuador
Set-Row 31 132475 1329 $null 18387 $null 23 11236
# Row 85 - Costa de Marfil
Set-Row 85 19556 55 19065 371 $null $null $null
# Row 102 - Guayana Francesa
Set-Row 102 9831 41 9472 294 $null $null $null
# Row 110 - Zimbabue
Set-Row 110 7787 35 6057 1503 $null $null $null
# Row 112 - Mauritania
Set-Row 112 7457 24 7070 226 $null $null $null
# Row 127 - Ruanda
Set-Row 127 4798 9 3080 1689 $null 2 29
# Row 135 - Aruba
Set-Row 135 3799 43 2753 1021 $null $null $null
# Row 144 - Mali
Set-Row 144 3064 23 2402 532 $null $null $null
# Row 159 - Togo
Set-Row 159 1722 15 1312 366 $null $null $null

# ---------------------------------------------------------------------
# 3) "Islas Malvinas" / "Montserrat" swap their relative order in the
#    country list; row 215 now reports Montserrat's figures and row 216
#    now reports Islas Malvinas' figures (full B:H rows trade places).
# ---------------------------------------------------------------------
$row215 = $ws.Range("A215:H215").Value2
$row216 = $ws.Range("A216:H216").Value2

$ws.Range("A216:H216").Value = $row215
$ws.Range("A215:H215").Value = $row216
